$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")
$lo = $ws.ListObjects.Item(1)

# Add a new row to the "Tabla24" table (this grows the table range and the
# autofilter automatically) for the new movie "Madame Web".
$newListRow = $lo.ListRows.Add()
$newRow = $newListRow.Range.Row

# The new row starts out unformatted; copy the formatting (alignment /
# number formats) from the row directly above it so it looks like every
# other row in the table before we fill in its values.
$ws.Range("B" + ($newRow - 1) + ":I" + ($newRow - 1)).Copy()
$ws.Range("B" + $newRow + ":I" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 2).Value = "Madame Web"
$ws.Cells.Item($newRow, 4).Value = 7
$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = 5
$ws.Cells.Item($newRow, 7).Value = 6
$ws.Cells.Item($newRow, 8).Value = 4.0999999999999996
$ws.Cells.Item($newRow, 9).Value = 3.7
$ws.Cells.Item($newRow, 3).Formula = "=AVERAGE(D" + $newRow + ",E" + $newRow + ",E" + $newRow + ",F" + $newRow + ",G" + $newRow + ",H" + $newRow + ",H" + $newRow + ",I" + $newRow + ")"

# Re-apply the table's descending sort on "Puntuación total" (column C) so
# the new row slots into its correct sorted position, like the rest of the
# table (this preserves every other row's values/formatting as-is).
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($lo.ListColumns.Item(2).DataBodyRange, 0, 2)
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Normalize the left-alignment style on column B for the two previously
# "newest" rows (Ironman / El Grinch) so they match the rest of the table
# now that "Madame Web" is the newest entry.
$ws.Range("B37").HorizontalAlignment = -4131
$ws.Range("B81").HorizontalAlignment = -4131

# Keep the sheet's active selection in sync with the bottom of the table.
$ws.Range("H124").Select()
